# Aggiorno file need_to_buy.xlsx da R
# The data window shifts forward by one day: row 2's old data is dropped,
# every remaining row's values are refreshed, and a new row (15) of data
# is appended for the newest day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45869
$ws.Range("B2").Value = 4654.8038038703
$ws.Range("C2").Value = 4546.33052745781
$ws.Range("D2").Value = 2376
$ws.Range("E2").Value = 5174.47883
$ws.Range("F2").Value = 112.083564732813

$ws.Range("A3").Value = 45870
$ws.Range("B3").Value = 5170.75192670309
$ws.Range("C3").Value = 4206.44960673379
$ws.Range("D3").Value = 1944
$ws.Range("E3").Value = 5482.523264
$ws.Range("F3").Value = 107.259206001279

$ws.Range("A4").Value = 45871
$ws.Range("B4").Value = 911.096850327142
$ws.Range("C4").Value = 1587.24321268832
$ws.Range("D4").Value = 1944
$ws.Range("E4").Value = 1753.382561
$ws.Range("F4").Value = 20.2303718067158

$ws.Range("A5").Value = 45872
$ws.Range("B5").Value = 803.704142719658
$ws.Range("C5").Value = 1566.57794126871
$ws.Range("D5").Value = 1944
$ws.Range("E5").Value = 1634.561035
$ws.Range("F5").Value = 18.8931180645439

$ws.Range("A6").Value = 45873
$ws.Range("B6").Value = 3848.98085240534
$ws.Range("C6").Value = 3818.44717252207
$ws.Range("D6").Value = 1944
$ws.Range("E6").Value = 4283.408683
$ws.Range("F6").Value = 96.2031251298634

$ws.Range("A7").Value = 45874
$ws.Range("B7").Value = 3939.76634278604
$ws.Range("C7").Value = 3810.14641676555
$ws.Range("D7").Value = 1944
$ws.Range("E7").Value = 4361.527973
$ws.Range("F7").Value = 95.3295019574796

$ws.Range("A8").Value = 45875
$ws.Range("B8").Value = 3939.76634278604
$ws.Range("C8").Value = 3809.85857930741
$ws.Range("D8").Value = 1944
$ws.Range("E8").Value = 4361.527973
$ws.Range("F8").Value = 95.3175087300571

$ws.Range("A9").Value = 45876
$ws.Range("B9").Value = 3939.76634278604
$ws.Range("C9").Value = 3764.59693853435
$ws.Range("D9").Value = 1944
$ws.Range("E9").Value = 4361.527973
$ws.Range("F9").Value = 93.4316070311794

$ws.Range("A10").Value = 45877
$ws.Range("B10").Value = 3939.76634278604
$ws.Range("C10").Value = 3662.75009455548
$ws.Range("D10").Value = 1944
$ws.Range("E10").Value = 4361.527973
$ws.Range("F10").Value = 89.1879885320599

$ws.Range("A11").Value = 45878
$ws.Range("B11").Value = 715.68594436081
$ws.Range("C11").Value = 1322.01866387421
$ws.Range("D11").Value = 1944
$ws.Range("E11").Value = 1499.552214
$ws.Range("F11").Value = 6.74520556305826

$ws.Range("A12").Value = 45879
$ws.Range("B12").Value = 634.559148236726
$ws.Range("C12").Value = 1290.42140457049
$ws.Range("D12").Value = 1944
$ws.Range("E12").Value = 1408.839072
$ws.Range("F12").Value = 5.02922201390678

$ws.Range("A13").Value = 45880
$ws.Range("B13").Value = 3735.57396241807
$ws.Range("C13").Value = 3564.89220001579
$ws.Range("D13").Value = 1944
$ws.Range("E13").Value = 4164.579533
$ws.Range("F13").Value = 85.4124071082382

$ws.Range("A14").Value = 45881
$ws.Range("B14").Value = 3735.57396241807
$ws.Range("C14").Value = 3598.90691384793
$ws.Range("D14").Value = 1944
$ws.Range("E14").Value = 4164.579533
$ws.Range("F14").Value = 86.829686851244

$ws.Range("A15").Value = 45882
$ws.Range("B15").Value = 3735.57396241807
$ws.Range("C15").Value = 3592.62224095091
$ws.Range("D15").Value = 1944
$ws.Range("E15").Value = 4164.579533
$ws.Range("F15").Value = 86.567825480535
